$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "azd-2_22-cv-02126"
$ws.Range("D2").Value = 100612
$ws.Range("E2").Value = 608

$ws.Range("A3").Value = "cand-3_18-cv-04865"
$ws.Range("D3").Value = 27697
$ws.Range("E3").Value = 615

$ws.Range("A4").Value = "cand-4_22-cv-02672"
$ws.Range("D4").Value = 44742
$ws.Range("E4").Value = 518

$ws.Range("A5").Value = "cand_22_cv_02094"
$ws.Range("D5").Value = 76610
$ws.Range("E5").Value = 458

$ws.Range("A6").Value = "cand_23_cv_02560"
$ws.Range("D6").Value = 67486
$ws.Range("E6").Value = 2157

$ws.Range("A7").Value = "cand_23_cv_03518"
$ws.Range("D7").Value = 32067
$ws.Range("E7").Value = 725

$ws.Range("A8").Value = "cand_24_cv_03170"
$ws.Range("D8").Value = 25448
$ws.Range("E8").Value = 980

$ws.Range("A9").Value = "cand_24_cv_04196"
$ws.Range("D9").Value = 11903
$ws.Range("E9").Value = 546

$ws.Range("A10").Value = "cand_3_22-cv-00956"
$ws.Range("D10").Value = 20671
$ws.Range("E10").Value = 421

$ws.Range("A11").Value = "casd_3_23-cv-01216"
$ws.Range("D11").Value = 34233
$ws.Range("E11").Value = 1002

$ws.Range("A12").Value = "ctd-3-23-cv-01035"
$ws.Range("D12").Value = 63448
$ws.Range("E12").Value = 585

$ws.Range("A13").Value = "dcd-1_23-cv-02055"
$ws.Range("D13").Value = 37155
$ws.Range("E13").Value = 476

$ws.Range("A14").Value = "dde_ 23_cv_1466"
$ws.Range("D14").Value = 34469
$ws.Range("E14").Value = 483

$ws.Range("A15").Value = "dde_21_cv_55"
$ws.Range("D15").Value = 45100
$ws.Range("E15").Value = 818

$ws.Range("A16").Value = "flsd-1_23-cv-23139"
$ws.Range("D16").Value = 16270
$ws.Range("E16").Value = 839

$ws.Range("A17").Value = "ilnd-1-21-cv-04349"
$ws.Range("D17").Value = 31039
$ws.Range("E17").Value = 726

$ws.Range("A18").Value = "mad-1-21-cv-10933"
$ws.Range("D18").Value = 16693
$ws.Range("E18").Value = 420

$ws.Range("A19").Value = "mied-4-23-cv-13132"
$ws.Range("D19").Value = 64927
$ws.Range("E19").Value = 556

$ws.Range("A20").Value = "nysd_20_cv_04494"
$ws.Range("D20").Value = 51071
$ws.Range("E20").Value = 953

$ws.Range("A21").Value = "nysd_22-cv-07111"
$ws.Range("D21").Value = 27485
$ws.Range("E21").Value = 744

$ws.Range("A22").Value = "nysd_22_cv_10292"
$ws.Range("D22").Value = 24365
$ws.Range("E22").Value = 448

$ws.Range("A23").Value = "nysd_23_cv_9476"
$ws.Range("D23").Value = 16202
$ws.Range("E23").Value = 422

$ws.Range("A24").Value = "nysd_24_cv_310"
$ws.Range("D24").Value = 43721
$ws.Range("E24").Value = 1634

$ws.Range("A25").Value = "txnd-4_24-cv-00673"
$ws.Range("D25").Value = 55525
$ws.Range("E25").Value = 584

$ws.Range("A26").Value = "txsd-4-21-cv-02473"
$ws.Range("D26").Value = 67106
$ws.Range("E26").Value = 492
